$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.440.79'
$ws.Range("E2").Value = '  -2.14%  '

$ws.Range("D3").Value = '2.891.75'
$ws.Range("E3").Value = '  -2.01%  '

$ws.Range("D4").Value = "'" + '0.998'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = "'" + '566.97'
$ws.Range("E5").Value = '  -4.45%  '

$ws.Range("D6").Value = "'" + '143.18'
$ws.Range("E6").Value = '  -3.25%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = "'" + '0.506'
$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("D9").Value = '2.891.28'
$ws.Range("E9").Value = '  -1.94%  '

$ws.Range("D10").Value = "'" + '6.81'
$ws.Range("E10").Value = '  -7.21%  '

$ws.Range("D11").Value = "'" + '0.146'
$ws.Range("E11").Value = '  -3.88%  '

$ws.Range("D12").Value = "'" + '0.434'
$ws.Range("E12").Value = '  -2.26%  '

$ws.Range("D13").Value = "'" + '0.0000235'
$ws.Range("E13").Value = '  -1.82%  '

$ws.Range("D14").Value = "'" + '31.90'
$ws.Range("E14").Value = '  -3.00%  '

$ws.Range("E15").Value = '  -0.62%  '

$ws.Range("D16").Value = '3.363.81'
$ws.Range("E16").Value = '  -2.22%  '

$ws.Range("D17").Value = '61.349.11'
$ws.Range("E17").Value = '  -2.28%  '

$ws.Range("D18").Value = "'" + '6.58'
$ws.Range("E18").Value = '  -2.16%  '

$ws.Range("D19").Value = '2.882.74'
$ws.Range("E19").Value = '  -2.24%  '

$ws.Range("D20").Value = "'" + '431.92'
$ws.Range("E20").Value = '  -2.43%  '

$ws.Range("D21").Value = "'" + '13.11'
$ws.Range("E21").Value = '  -2.76%  '

$ws.Range("E22").Value = '  -2.09%  '

$ws.Range("E23").Value = '  -3.10%  '

$ws.Range("D24").Value = "'" + '79.09'
$ws.Range("E24").Value = '  -2.89%  '

$ws.Range("D25").Value = "'" + '11.84'
$ws.Range("E25").Value = '  +1.31%  '

$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = "'" + '10.02'
$ws.Range("E26").Value = '  -10.14%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = "'" + '1.00'
$ws.Range("E27").Value = '  +0.05%  '

$ws.Range("D28").Value = "'" + '2.01'
$ws.Range("E28").Value = '  -6.78%  '

$ws.Range("D29").Value = "'" + '0.0000104'
$ws.Range("E29").Value = '  -1.40%  '

$ws.Range("D30").Value = "'" + '7.00'
$ws.Range("E30").Value = '  -3.67%  '

$ws.Range("E31").Value = '  -4.80%  '

$ws.Range("E32").Value = '  -8.34%  '

$ws.Range("D33").Value = "'" + '0.998'
$ws.Range("E33").Value = '  -0.29%  '

$ws.Range("E34").Value = '  -2.51%  '

$ws.Range("D35").Value = "'" + '25.55'
$ws.Range("E35").Value = '  -3.71%  '

$ws.Range("D36").Value = "'" + '0.954'
$ws.Range("E36").Value = '  -3.90%  '

$ws.Range("E37").Value = '  -3.77%  '

$ws.Range("D38").Value = "'" + '48.84'
$ws.Range("E38").Value = '  -1.64%  '

$ws.Range("E39").Value = '  -5.31%  '

$ws.Range("E40").Value = '  -12.04%  '

$ws.Range("E41").Value = '  -3.28%  '

$ws.Range("E42").Value = '  -3.59%  '

$ws.Range("D43").Value = "'" + '39.49'
$ws.Range("E43").Value = '  -1.04%  '

$ws.Range("D44").Value = "'" + '0.268'
$ws.Range("E44").Value = '  -5.04%  '

$ws.Range("D45").Value = '2.687.14'
$ws.Range("E45").Value = '  -0.52%  '

$ws.Range("D46").Value = "'" + '133.60'
$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("E47").Value = '  -1.61%  '

$ws.Range("D49").Value = "'" + '340.09'
$ws.Range("E49").Value = '  -6.79%  '

$ws.Range("E50").Value = '  -1.88%  '

$ws.Range("D51").Value = "'" + '21.46'
$ws.Range("E51").Value = '  -6.59%  '
